$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; existing rows 1-5 shift down to 2-6.
$ws.Rows.Item(1).Insert()

# --- New row 1 content ---
# A1:C1 merged -> "고객명" (customer name) label, centered, bottom border only
$ws.Range("A1:C1").Merge()
$ws.Range("A1").Value = "고객명"

# N1 -> "문서생성날짜" (document creation date) label, centered, no border
$ws.Range("N1").Value = "문서생성날짜"

# O1:P1 merged -> "날짜" value area, centered, bottom border only
$ws.Range("O1:P1").Merge()
$ws.Range("O1").Value = "날짜"

# Styling for the label / value cells of the new row
$labelRange = $ws.Range("A1:C1,O1:P1")
$labelRange.HorizontalAlignment = -4108
$labelRange.VerticalAlignment = -4108
$labelRange.Borders.Item(9).LineStyle = 1

$n1 = $ws.Range("N1")
$n1.HorizontalAlignment = -4108
$n1.VerticalAlignment = -4108

# --- Selection / active cell to match saved view ---
$ws.Range("A2:A4").Select()
